$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 9434
$ws1.Range("F6").Value = 9434
$ws1.Range("F7").Value = 606
$ws1.Range("F19").Value = 89
$ws1.Range("F33").Value = 70
$ws1.Range("F36").Value = 2146
$ws1.Range("F39").Value = 3637
$ws1.Range("F40").Value = 544
$ws1.Range("F41").Value = 2625
$ws1.Range("F42").Value = 3056
$ws1.Range("F43").Value = 1323
$ws1.Range("F47").Value = 526
$ws1.Range("F48").Value = 68
$ws1.Range("F49").Value = 226

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 9434
$ws4.Range("F10").Value = 606
$ws4.Range("F34").Value = 70
$ws4.Range("F38").Value = 2146
$ws4.Range("F42").Value = 3637
$ws4.Range("F43").Value = 3056
$ws4.Range("F45").Value = 1323
$ws4.Range("F49").Value = 526
$ws4.Range("F50").Value = 68
$ws4.Range("F51").Value = 226
